$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 182, shifting existing rows 182:192 down to 183:193
$ws.Rows(182).Insert()

# Populate the newly inserted row 182 with the new weekly price record
$ws.Range("A182").Value = 8
$ws.Range("B182").Value = "Terminal La Palmera de La Serena"
$ws.Range("C182").Value = "Coquimbo"
$ws.Range("D182").Value = 45223
$ws.Range("E182").Value = 4
$ws.Range("F182").Value = 100114007
$ws.Range("G182").Value = "Jengibre"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 460
$ws.Range("K182").Value = 28000
$ws.Range("L182").Value = 29000
$ws.Range("M182").Value = 28500
$ws.Range("N182").Value = "$/caja 13 kilos"
$ws.Range("O182").Value = "Perú"
$ws.Range("P182").Value = 2192
$ws.Range("Q182").Value = 13
$ws.Range("R182").Value = "Hortaliza"
